$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana 1")
$ws.Range("D16").Select()
$excel.ActiveWindow.ScrollRow = 2
